# "incorporate negaoctet production emissions"
# Updates the PC emissions model: refreshed power-draw/usage-life inputs for
# Desktop/Laptop/Monitor and a higher-precision (4-decimal) number format for
# the PEPS column now that the inputs carry more signal.

$wb = $excel.ActiveWorkbook
$pc = $wb.Worksheets.Item("PC")

# --- Updated source inputs (Desktop / Laptop / Monitor rows) ---------------
# Installed base (M) -- now literal figures rather than the old min/max avg
$pc.Range("F10").Value = 277
$pc.Range("F11").Value = 175
$pc.Range("F12").Value = 69

# UL (years)
$pc.Range("J10").Value = 6
$pc.Range("J11").Value = 5

# --- Number formatting ------------------------------------------------------
# Bump the PEPS (gCO2e/s) column to a 4-decimal accounting format so the
# smaller recomputed values remain legible.
$pc.Range("K10:K15").NumberFormat = '_(* #,##0.0000_);_(* \(#,##0.0000\);_(* "-"??_);_(@_)'

# Keep the surrounding input columns on their existing single-decimal format.
$pc.Range("F10:F14").NumberFormat = "0.0"
$pc.Range("I10:I12").NumberFormat = "0.0"
$pc.Range("I14").NumberFormat = "0.0"
$pc.Range("J10:J12").NumberFormat = "0.0"

# --- View state: PC becomes the active tab/selection -----------------------
$pc.Activate()
$pc.Range("J13").Select()
